$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$ws.Range("G61").Select()
$panes = $win.Panes
$panes.Item(2).ScrollRow = 46
$panes.Item(2).ScrollColumn = 1
Write-Host "Pane2 ScrollRow=$($panes.Item(2).ScrollRow) ScrollColumn=$($panes.Item(2).ScrollColumn)"
